$wb = $excel.ActiveWorkbook

# Rename existing Sheet1 -> Login
$wsLogin = $wb.Worksheets.Item(1)
$wsLogin.Name = "Login"

# Add Search sheet right after Login
$wsSearch = $wb.Worksheets.Add($null, $wsLogin)
$wsSearch.Name = "Search"

# Add new blank Sheet1 right after Search
$wsSheet1 = $wb.Worksheets.Add($null, $wsSearch)
$wsSheet1.Name = "Sheet1"

# Populate Search sheet with category/subcategory test data
$data = @(
    @("category", "subcategory"),
    @("WOMEN", "DRESS"),
    @("WOMEN", "TOPS"),
    @("WOMEN", "SAREE"),
    @("MEN", "TSHIRTS"),
    @("MEN", "JEANS"),
    @("KIDS", "DRESS"),
    @("KIDS", "TOPS & SHIRTS")
)
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $wsSearch.Cells.Item($row, 1).Value = $data[$i][0]
    $wsSearch.Cells.Item($row, 2).Value = $data[$i][1]
}

# Best-fit width for column B (subcategory) on the Search sheet
$wsSearch.Columns.Item(2).ColumnWidth = 12.25

# Restore each sheet's last-used selection
$wsLogin.Range("B38").Select() | Out-Null
$wsSearch.Range("D9").Select() | Out-Null
$wsSheet1.Range("A1:B6").Select() | Out-Null

# Search is the active (visible) tab
$wsSearch.Activate() | Out-Null
